$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Numeric cell updates ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 67
$ws.Range("K16").Value = -14.925373134328
$ws.Range("L16").Value = 18.75
$ws.Range("M16").Value = -3.389830508474
$ws.Range("N16").Value = -83.185840707964
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -85.714285714285
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -75
$ws.Range("I17").Value = 72
$ws.Range("J17").Value = 82
$ws.Range("K17").Value = -12.195121951219
$ws.Range("L17").Value = 22.033898305084
$ws.Range("M17").Value = 84.615384615384
$ws.Range("N17").Value = -20
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 12
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 4.166666666666
$ws.Range("I18").Value = 85
$ws.Range("J18").Value = 98
$ws.Range("K18").Value = -13.265306122449
$ws.Range("L18").Value = 21.428571428571
$ws.Range("M18").Value = -19.047619047619
$ws.Range("N18").Value = -84.739676840215
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 20
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -12.987012987013
$ws.Range("I19").Value = 336
$ws.Range("J19").Value = 318
$ws.Range("K19").Value = 5.66037735849
$ws.Range("L19").Value = 42.978723404255
$ws.Range("M19").Value = -32.121212121212
$ws.Range("N19").Value = -62.331838565022
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 22
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = 29.411764705882
$ws.Range("L20").Value = 214.285714285714
$ws.Range("M20").Value = 120
$ws.Range("N20").Value = -94.965675057208
$ws.Range("C21").Value = 31
$ws.Range("E21").Value = -20.51282051282
$ws.Range("F21").Value = 115
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -21.232876712328
$ws.Range("I21").Value = 576
$ws.Range("J21").Value = 591
$ws.Range("K21").Value = -2.538071065989
$ws.Range("L21").Value = 35.211267605633
$ws.Range("M21").Value = -19.214586255259
$ws.Range("N21").Value = -75.193798449612
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 39
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = 85.714285714285
$ws.Range("M22").Value = -10.344827586206
$ws.Range("M23").Value = -50
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 58
$ws.Range("E24").Value = -27.586206896551
$ws.Range("F24").Value = 182
$ws.Range("G24").Value = 202
$ws.Range("H24").Value = -9.900990099009
$ws.Range("I24").Value = 699
$ws.Range("J24").Value = 810
$ws.Range("K24").Value = -13.703703703703
$ws.Range("L24").Value = 15.346534653465
$ws.Range("M24").Value = 16.30615640599
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = -36.363636363636
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -9.756097560975
$ws.Range("I25").Value = 162
$ws.Range("J25").Value = 173
$ws.Range("K25").Value = -6.35838150289
$ws.Range("L25").Value = 43.362831858407
$ws.Range("M25").Value = 8
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -33.333333333333
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 57.142857142857
$ws.Range("I27").Value = 43
$ws.Range("J27").Value = 32
$ws.Range("K27").Value = 34.375
$ws.Range("L27").Value = 43.333333333333

# --- Cells that become text placeholders ("0" / "***.*") ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C30").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "***.*"
$ws.Range("E30").NumberFormat = "General"
